# Update countries & provincias Spain
# - Reorders a few country rows (swap adjacent country-name cells)
# - Refreshes the "last updated" timestamp
# - Refreshes the COVID case numbers for a set of countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 16:06"

# --- Swap country-name pairs (Irak now appears before Emiratos Arabes Unidos, etc.) ---
$ws.Range("A34").Value = "Irak"
$ws.Range("A35").Value = "Emiratos Arabes Unidos"

$ws.Range("A80").Value = "Kenia"
$ws.Range("A81").Value = "Republica de Macedonia"

$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"

$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Refresh case numbers (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos
$ws.Range("B4").Value = 2683902
$ws.Range("C4").Value = 2091
$ws.Range("D4").Value = 1122593
$ws.Range("E4").Value = 1432467
$ws.Range("G4").Value = 59
$ws.Range("H4").Value = 128842

# India
$ws.Range("B7").Value = 572723
$ws.Range("C7").Value = 5187
$ws.Range("D7").Value = 338516
$ws.Range("E7").Value = 217222
$ws.Range("G7").Value = 81
$ws.Range("H7").Value = 16985

# Arabia Saudita
$ws.Range("B18").Value = 190823
$ws.Range("C18").Value = 4387
$ws.Range("D18").Value = 130766
$ws.Range("E18").Value = 58408
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = 1649

# Catar
$ws.Range("B23").Value = 96088
$ws.Range("C23").Value = 982
$ws.Range("D23").Value = 81564
$ws.Range("E23").Value = 14411

# Emiratos Arabes Unidos (row 34, after the swap above)
$ws.Range("B34").Value = 49109
$ws.Range("C34").Value = 1958
$ws.Range("D34").Value = 24760
$ws.Range("E34").Value = 22406
$ws.Range("G34").Value = 104
$ws.Range("H34").Value = 1943

# Irak (row 35, after the swap above)
$ws.Range("B35").Value = 48667
$ws.Range("C35").Value = 421
$ws.Range("D35").Value = 37566
$ws.Range("E35").Value = 10786
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 315

# Portugal
$ws.Range("B39").Value = 42141
$ws.Range("C39").Value = 229
$ws.Range("D39").Value = 27505
$ws.Range("E39").Value = 13060
$ws.Range("G39").Value = 8
$ws.Range("H39").Value = 1576

# Suiza
$ws.Range("E46").Value = 651
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 1963

# Serbia
$ws.Range("B62").Value = 14564
$ws.Range("C62").Value = 276
$ws.Range("D62").Value = 12662
$ws.Range("E62").Value = 1625
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 277

# Kenia (row 80, after the swap above)
$ws.Range("B80").Value = 6366
$ws.Range("C80").Value = 176
$ws.Range("D80").Value = 2013
$ws.Range("E80").Value = 4209
$ws.Range("H80").Value = 144

# Republica de Macedonia (row 81, after the swap above)
$ws.Range("B81").Value = 6334
$ws.Range("C81").Value = 125
$ws.Range("D81").Value = 2475
$ws.Range("E81").Value = 3557
$ws.Range("G81").Value = 4
$ws.Range("H81").Value = 302

# Guayana Francesa
$ws.Range("B95").Value = 4004
$ws.Range("C95").Value = 230
$ws.Range("D95").Value = 1508
$ws.Range("E95").Value = 2481

# Mayotte
$ws.Range("B102").Value = 2603
$ws.Range("C102").Value = 43
$ws.Range("E102").Value = 243
$ws.Range("G102").Value = 3
$ws.Range("H102").Value = 35

# Estado de Palestina
$ws.Range("E104").Value = 1933
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 6

# Cuba
$ws.Range("B105").Value = 2341
$ws.Range("C105").Value = 1
$ws.Range("D105").Value = 2214
$ws.Range("E105").Value = 41

# Sri Lanka
$ws.Range("B111").Value = 2047
$ws.Range("C111").Value = 8
$ws.Range("E111").Value = 325

# Libano
$ws.Range("B117").Value = 1778
$ws.Range("C117").Value = 33
$ws.Range("D117").Value = 1183
$ws.Range("E117").Value = 561

# Liberia
$ws.Range("B145").Value = 780
$ws.Range("C145").Value = 10
$ws.Range("D145").Value = 324
$ws.Range("E145").Value = 420

# Reunion
$ws.Range("B153").Value = 526
$ws.Range("C153").Value = 4
$ws.Range("E153").Value = 52
